# Apply cryptos list update (Sat Apr 29 04:09:55 UTC 2023 snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.493.25"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.907.61"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'326.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'0.4844"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("D8").Value = "'0.4072"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "'0.08135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").Value = "'23.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("D12").Value = "1.917.47"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "'6.022"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'7.096"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'90.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'0.06757"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "'0.00001042"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").Value = "'17.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "29.514.23"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'5.607"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'2.165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "2.129.01"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'155.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "'20.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "'6.293"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.93%  "
$ws.Range("D29").Value = "'2.103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "'119.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").Value = "'1.036"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "'0.09562"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'5.536"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").Value = "'1.394"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'3.551"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02265"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'1.172"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'0.5959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.921"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'10.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("D42").Value = "'0.1855"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "'2.434"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").Value = "'0.07721"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'12.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("D47").Value = "'0.5575"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").Value = "'1.957"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "'115.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'72.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  +2.37%  "
